$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.560.52"
$ws.Range("E2").Value = "  +1.97%  "

# Row 3
$ws.Range("D3").Value = "3.095.15"
$ws.Range("E3").Value = "  +0.37%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'525.09"
$ws.Range("E5").Value = "  +1.98%  "

# Row 6
$ws.Range("E6").Value = "  +1.79%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.442"
$ws.Range("E8").Value = "  +1.92%  "

# Row 9
$ws.Range("E9").Value = "  +1.24%  "

# Row 10
$ws.Range("E10").Value = "  +0.36%  "

# Row 11
$ws.Range("E11").Value = "  +3.08%  "

# Row 12
$ws.Range("D12").Value = "3.624.27"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("E13").Value = "  +0.91%  "

# Row 14
$ws.Range("D14").Value = "'26.88"
$ws.Range("E14").Value = "  +4.84%  "

# Row 15
$ws.Range("E15").Value = "  +1.53%  "

# Row 16
$ws.Range("D16").Value = "58.583.78"
$ws.Range("E16").Value = "  +1.82%  "

# Row 17
$ws.Range("D17").Value = "3.091.16"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18
$ws.Range("E18").Value = "  +0.41%  "

# Row 19
$ws.Range("D19").Value = "'12.91"
$ws.Range("E19").Value = "  -1.31%  "

# Row 20
$ws.Range("D20").Value = "'8.11"
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("D21").Value = "'340.84"
$ws.Range("E21").Value = "  +1.77%  "

# Row 22
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("E23").Value = "  +0.82%  "

# Row 24
$ws.Range("D24").Value = "'65.86"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("E25").Value = "  +0.13%  "

# Row 26
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("E27").Value = "  -0.93%  "

# Row 28
$ws.Range("D28").Value = "'6.60"
$ws.Range("E28").Value = "  +2.76%  "

# Row 29
$ws.Range("D29").Value = "'7.24"
$ws.Range("E29").Value = "  +1.97%  "

# Row 30
$ws.Range("D30").Value = "'1.86"
$ws.Range("E30").Value = "  +2.64%  "

# Row 31
$ws.Range("D31").Value = "'21.02"
$ws.Range("E31").Value = "  +0.77%  "

# Row 32
$ws.Range("E32").Value = "  +3.20%  "

# Row 33
$ws.Range("D33").Value = "'154.28"
$ws.Range("E33").Value = "  +0.26%  "

# Row 34
$ws.Range("D34").Value = "'4.62"
$ws.Range("E34").Value = "  +1.86%  "

# Row 35
$ws.Range("D35").Value = "'6.07"
$ws.Range("E35").Value = "  +2.89%  "

# Row 36
$ws.Range("D36").Value = "'26.92"
$ws.Range("E36").Value = "  -1.46%  "

# Row 37
$ws.Range("D37").Value = "'1.32"
$ws.Range("E37").Value = "  +6.76%  "

# Row 38
$ws.Range("E38").Value = "  -0.53%  "

# Row 39
$ws.Range("D39").Value = "3.135.59"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("E40").Value = "  +1.30%  "

# Row 41
$ws.Range("D41").Value = "'36.82"
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.671"
$ws.Range("E42").Value = "  -0.15%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.48"
$ws.Range("E44").Value = "  +6.45%  "

# Row 45
$ws.Range("D45").Value = "2.277.91"
$ws.Range("E45").Value = "  -0.50%  "

# Row 46
$ws.Range("D46").Value = "'0.0257"
$ws.Range("E46").Value = "  +1.44%  "

# Row 47
$ws.Range("D47").Value = "'20.85"
$ws.Range("E47").Value = "  +3.71%  "

# Row 48
$ws.Range("D48").Value = "'0.961"
$ws.Range("E48").Value = "  +2.06%  "

# Row 49
$ws.Range("D49").Value = "'6.01"
$ws.Range("E49").Value = "  +2.09%  "

# Row 50
$ws.Range("D50").Value = "'267.46"
$ws.Range("E50").Value = "  +8.36%  "

# Row 51
$ws.Range("D51").Value = "'0.747"
$ws.Range("E51").Value = "  +8.75%  "
